$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point ("2026/01/11", time 13) arrived, so a fresh row is
# inserted right before the current row 601, pushing the existing
# 601:642 block down to 602:643 (dimension grows from D642 to D643).
$ws.Rows.Item(601).Insert()

# Column A holds plain text dates (not real Excel dates) in this sheet,
# so force text formatting while assigning the value to stop Excel from
# auto-parsing "2026/01/11" into a date serial number, then restore the
# default "Normal" style so the new row matches its neighbours.
$ws.Cells.Item(601, 1).NumberFormat = "@"
$ws.Cells.Item(601, 1).Value = "2026/01/11"
$ws.Cells.Item(601, 1).Style = "Normal"

$ws.Cells.Item(601, 2).Value = "日"
$ws.Cells.Item(601, 3).Value = 13
$ws.Cells.Item(601, 4).Value = 12
